$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" column header in H1 — copy formatting (bold, border, centered/top
# alignment) from the existing header cell G1 so it matches the rest of the
# header row, then set its text.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = $false
$ws.Range("H1").Value = "Save"

# Corresponding data value for the new column
$ws.Range("H2").Value = 0
